$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the Area / Atotal columns
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Side-by-side comparison columns (J: Atotal, K: Qtotal)
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Row 2: area of first segment, running total, and the comparison cells
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Row 3 area (not part of the G4:G15 shared-formula block)
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15 share one formula, mirroring the existing D/E shared formulas
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Put the selection where the author left it
$ws.Range("J2:K2").Select() | Out-Null
